$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily rows appended after the existing data (row 52 was the last one).
$newRows = @(
    @(46010, 5600, 3952, 3629, 209, 49, 59, 5, 1),
    @(46013, 5590, 3758, 3459, 214, 50, 31, 3, 1),
    @(46014, 5591, 3624, 3325, 201, 48, 46, 4, 0)
)

$startRow = 53
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $vals = $newRows[$i]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
    # Column A carries the date number format used by the other date cells.
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($startRow - 1, 1).NumberFormat
}

$lastRow = $startRow + $newRows.Length - 1
$ws.Range("A$($lastRow):I$($lastRow)").Select()
$excel.ActiveWindow.ScrollRow = 43
